# "control viajes modificacion completa"
# Rebuild the "Producto N / Cantidad N" pairs for rows 4, 7, 8 and tweak a
# few summary columns on row 9, per the pedidos_calendario.xlsx update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 ----------------------------------------------------------
$ws.Range("E4").Value = $false
$ws.Range("G4").Value = 0.0
$ws.Range("H4").Value = 100.0
$ws.Range("I4").Value = 1.0
$ws.Range("J4").Value = 60.0
$ws.Range("K4").Value = 4.0
$ws.Range("L4").Value = 40.0

# ---- Row 7 ----------------------------------------------------------
$ws.Range("E7").Value = $false
$ws.Range("G7").Value = 0.0
$ws.Range("H7").Value = 150.0
$ws.Range("I7").Value = 1.0
$ws.Range("J7").Value = 1150.0
$ws.Range("K7").Value = 2.0
$ws.Range("L7").Value = 1515.0
$ws.Range("M7").Value = 3.0
$ws.Range("N7").Value = 10150.0
$ws.Range("O7").Value = 4.0
$ws.Range("P7").Value = 150.0
$ws.Range("Q7").Value = 5.0
$ws.Range("R7").Value = 1150.0
$ws.Range("S7").Value = 6.0
$ws.Range("T7").Value = 150.0
$ws.Range("U7").Value = 7.0
$ws.Range("V7").Value = 115.0
$ws.Range("W7").Value = 8.0
$ws.Range("X7").Value = 150.0

# ---- Row 8 ----------------------------------------------------------
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $false
$ws.Range("H8").Value = 132.0
$ws.Range("J8").Value = 240.0
$ws.Range("L8").Value = 1050.0
$ws.Range("N8").Value = 300.0
$ws.Range("P8").Value = 410.0
$ws.Range("R8").Value = 310.0
$ws.Range("T8").Value = 164.0
$ws.Range("V8").Value = 200.0
$ws.Range("W8").Value = 8.0
$ws.Range("X8").Value = 300.0

# ---- Row 9 ----------------------------------------------------------
$ws.Range("C9").Value = 1.0
$ws.Range("D9").Value = 7.0
$ws.Range("E9").Value = $false
